# smartfridge_budget.xlsx update
# Project changes: pivot towards interconnection with smartwatch-type connected objects.
# - Translate headers & category/description text to French
# - Update a few cost values
# - Insert a new "Marge imprévus (10%)" row before TOTAL
# - Update TOTAL accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 so the existing TOTAL row (row 8) shifts down to row 9.
$ws.Rows.Item(8).Insert()

# --- Header row ---
$ws.Range("A1").Value = "Catégorie"
$ws.Range("B1").Value = "Détail"
$ws.Range("C1").Value = "Coût estimé (€)"

# --- Row 2 ---
$ws.Range("A2").Value = "Prototype matériel"
$ws.Range("B2").Value = "Frigo modifié + caméras + module IA embarqué"
$ws.Range("C2").Value = 1100

# --- Row 3 ---
$ws.Range("A3").Value = "Développement logiciel"
$ws.Range("B3").Value = "App mobile, backend FastAPI, intégration Drive"
$ws.Range("C3").Value = 0

# --- Row 4 ---
$ws.Range("A4").Value = "Ressources humaines (IA / IoT / Mobile / Backend)"
$ws.Range("B4").Value = "1 mois de travail cumulé des rôles clés"
$ws.Range("C4").Value = 4200

# --- Row 5 ---
$ws.Range("A5").Value = "Infra Cloud (hébergement API / DB)"
$ws.Range("B5").Value = "Serveur OVH/AWS + stockage données 1 an"
$ws.Range("C5").Value = 500

# --- Row 6 ---
$ws.Range("A6").Value = "Tests terrain / panels utilisateurs"
$ws.Range("B6").Value = "Séances de tests en conditions réelles"
$ws.Range("C6").Value = 300

# --- Row 7 ---
$ws.Range("A7").Value = "Communication / gouvernance projet"
$ws.Range("B7").Value = "Réunions projet, doc, pilotage PO"
$ws.Range("C7").Value = 100

# --- Row 8 (new row, inserted above) ---
$ws.Range("A8").Value = "Marge imprévus (10%)"
$ws.Range("B8").Value = "Buffer sécurité budget matériel / délai fournisseur"
$ws.Range("C8").Value = 620

# --- Row 9 (former TOTAL row, now shifted down) ---
# B9 is intentionally left untouched: the row insert already carried the
# previously-empty B8 cell down to B9, keeping it an empty string cell.
$ws.Range("A9").Value = "TOTAL"
$ws.Range("C9").Value = 6820
